# "in te leveren lijst aangepast" - clean up the deliverables checklist:
#  - drop the duplicate "Week 3.2" sub-bullet under "Agenda's" (right after "Week 3.1")
#  - drop the duplicate "Week 3.2" sub-bullet under "Notules" (right after "Week 2.1")
#  - drop the stray "Week 3" sub-bullet under "Weekverslag"
#  - move the (cursor-position) "_GoBack" bookmark from "Weekverslag" to "Gedocumenteerde code"

$d = $word.ActiveDocument

# Remove the duplicated "Week 3.2" bullet right after "Week 3.1".
$d.Paragraphs.Item(3).Range.Delete()

# Remove the duplicated "Week 3.2" bullet right after "Week 2.1" (indices shifted
# by one now that the previous paragraph is gone).
$d.Paragraphs.Item(5).Range.Delete()

# Remove the stray "Week 3" bullet right after "Weekverslag".
$d.Paragraphs.Item(7).Range.Delete()

# Relocate the "_GoBack" bookmark from the "Weekverslag" paragraph onto the
# start of the "Gedocumenteerde code" paragraph.
$d.Bookmarks.Item("_GoBack").Delete()

$target = $d.Paragraphs.Item(7)
$bmRange = $target.Range.Duplicate
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)
